$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    4224.603527133314,
    4110.151735092823,
    4103.94636274034,
    4103.94636274034,
    4103.94636274034,
    4103.94636274034,
    4103.94636274034,
    4103.94636274034,
    4103.94636274034,
    3962.874990755126,
    3900.476813571693
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $values[$i]
}
